$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 2).Value = 6
$ws.Cells.Item(8, 3).Value = 'Alex Caruso'
$ws.Cells.Item(8, 4).Value = 'PG'
$ws.Cells.Item(8, 5).Value = '6-4'
$ws.Cells.Item(8, 6).Value = 186
$ws.Cells.Item(8, 7).Value = 'February 28, 1994'
$ws.Cells.Item(8, 8).Value = 'us'
$ws.Cells.Item(8, 9).Value = '5'
$ws.Cells.Item(8, 10).Value = 'Texas A&M'
$ws.Cells.Item(8, 11).Value = 'https://www.basketball-reference.com/players/c/carusal01.html'

# Row 9
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = 'Andre Drummond'
$ws.Cells.Item(9, 4).Value = 'C'
$ws.Cells.Item(9, 5).Value = '6-10'
$ws.Cells.Item(9, 6).Value = 279
$ws.Cells.Item(9, 7).Value = 'August 10, 1993'
$ws.Cells.Item(9, 8).Value = 'us'
$ws.Cells.Item(9, 9).Value = '10'
$ws.Cells.Item(9, 10).Value = 'UConn'
$ws.Cells.Item(9, 11).Value = 'https://www.basketball-reference.com/players/d/drumman01.html'

# Row 10
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = 'Derrick Jones Jr.'
$ws.Cells.Item(10, 4).Value = 'SF'
$ws.Cells.Item(10, 5).Value = '6-5'
$ws.Cells.Item(10, 6).Value = 210
$ws.Cells.Item(10, 7).Value = 'February 15, 1997'
$ws.Cells.Item(10, 8).Value = 'us'
$ws.Cells.Item(10, 9).Value = '6'
$ws.Cells.Item(10, 10).Value = 'UNLV'
$ws.Cells.Item(10, 11).Value = 'https://www.basketball-reference.com/players/j/jonesde02.html'

# Row 11
$ws.Cells.Item(11, 2).Value = 24
$ws.Cells.Item(11, 3).Value = 'Javonte Green'
$ws.Cells.Item(11, 4).Value = 'SG'
$ws.Cells.Item(11, 5).Value = '6-4'
$ws.Cells.Item(11, 6).Value = 205
$ws.Cells.Item(11, 7).Value = 'July 23, 1993'
$ws.Cells.Item(11, 8).Value = 'us'
$ws.Cells.Item(11, 9).Value = '3'
$ws.Cells.Item(11, 10).Value = 'Radford'
$ws.Cells.Item(11, 11).Value = 'https://www.basketball-reference.com/players/g/greenja02.html'

# Row 12
$ws.Cells.Item(12, 2).Value = 25
$ws.Cells.Item(12, 3).Value = 'Dalen Terry'
$ws.Cells.Item(12, 4).Value = 'SG'
$ws.Cells.Item(12, 5).Value = '6-7'
$ws.Cells.Item(12, 6).Value = 195
$ws.Cells.Item(12, 7).Value = 'July 12, 2002'
$ws.Cells.Item(12, 8).Value = 'us'
$ws.Cells.Item(12, 9).Value = 'R'
$ws.Cells.Item(12, 10).Value = 'Arizona'
$ws.Cells.Item(12, 11).Value = 'https://www.basketball-reference.com/players/t/terryda01.html'

# Row 13
$ws.Cells.Item(13, 2).Value = 21
$ws.Cells.Item(13, 3).Value = 'Patrick Beverley'
$ws.Cells.Item(13, 4).Value = 'PG'
$ws.Cells.Item(13, 5).Value = '6-1'
$ws.Cells.Item(13, 6).Value = 180
$ws.Cells.Item(13, 7).Value = 'July 12, 1988'
$ws.Cells.Item(13, 8).Value = 'us'
$ws.Cells.Item(13, 9).Value = '10'
$ws.Cells.Item(13, 10).Value = 'Arkansas'
$ws.Cells.Item(13, 11).Value = 'https://www.basketball-reference.com/players/b/beverpa01.html'

# Row 14
$ws.Cells.Item(14, 2).Value = 19
$ws.Cells.Item(14, 3).Value = 'Marko Simonovic'
$ws.Cells.Item(14, 4).Value = 'C'
$ws.Cells.Item(14, 5).Value = '6-11'
$ws.Cells.Item(14, 6).Value = 216
$ws.Cells.Item(14, 7).Value = 'October 15, 1999'
$ws.Cells.Item(14, 8).Value = 'me'
$ws.Cells.Item(14, 9).Value = '1'
$ws.Cells.Item(14, 10).ClearContents()
$ws.Cells.Item(14, 11).Value = 'https://www.basketball-reference.com/players/s/simonma01.html'

# Row 15
$ws.Cells.Item(15, 2).Value = 32
$ws.Cells.Item(15, 3).Value = 'Terry Taylor (TW)'
$ws.Cells.Item(15, 4).Value = 'PF'
$ws.Cells.Item(15, 5).Value = '6-5'
$ws.Cells.Item(15, 6).Value = 230
$ws.Cells.Item(15, 7).Value = 'September 23, 1999'
$ws.Cells.Item(15, 8).Value = 'us'
$ws.Cells.Item(15, 9).Value = '1'
$ws.Cells.Item(15, 10).Value = 'Austin Peay State University'
$ws.Cells.Item(15, 11).Value = 'https://www.basketball-reference.com/players/t/taylote01.html'

# Row 16
$ws.Cells.Item(16, 2).Value = 22
$ws.Cells.Item(16, 3).Value = 'Carlik Jones'
$ws.Cells.Item(16, 4).Value = 'PG'
$ws.Cells.Item(16, 5).Value = '6-1'
$ws.Cells.Item(16, 6).Value = 185
$ws.Cells.Item(16, 7).Value = 'December 23, 1997'
$ws.Cells.Item(16, 8).Value = 'us'
$ws.Cells.Item(16, 9).Value = '1'
$ws.Cells.Item(16, 10).Value = 'Radford, Louisville'
$ws.Cells.Item(16, 11).Value = 'https://www.basketball-reference.com/players/j/jonesca03.html'

# Row 17
$ws.Cells.Item(17, 2).ClearContents()
$ws.Cells.Item(17, 3).Value = 'Lonzo Ball'
$ws.Cells.Item(17, 4).Value = 'PG'
$ws.Cells.Item(17, 5).Value = '6-6'
$ws.Cells.Item(17, 6).Value = 190
$ws.Cells.Item(17, 7).Value = 'October 27, 1997'
$ws.Cells.Item(17, 8).Value = 'us'
$ws.Cells.Item(17, 9).Value = '5'
$ws.Cells.Item(17, 10).Value = 'UCLA'
$ws.Cells.Item(17, 11).Value = 'https://www.basketball-reference.com/players/b/balllo01.html'

# Row 18
$ws.Cells.Item(18, 2).ClearContents()
$ws.Cells.Item(18, 3).Value = 'Justin Lewis (TW)'
$ws.Cells.Item(18, 4).Value = 'PF'
$ws.Cells.Item(18, 5).Value = '6-7'
$ws.Cells.Item(18, 6).Value = 245
$ws.Cells.Item(18, 7).Value = 'April 12, 2002'
$ws.Cells.Item(18, 8).Value = 'us'
$ws.Cells.Item(18, 9).Value = 'R'
$ws.Cells.Item(18, 10).Value = 'Marquette'
$ws.Cells.Item(18, 11).Value = 'https://www.basketball-reference.com/players/l/lewisju02.html'
